# Update the BOM workbook: PPCB (Per Unit) and Case (Per Box) tables.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("D1").Value = "Per Box"

# --- Column headers (row 2) -------------------------------------------
$ws.Range("A2").Value = "Component"
$ws.Range("B2").Value = "Price"
$ws.Range("D2").Value = "Component"
$ws.Range("E2").Value = "Price"

# --- Row 3: Seeed XIAO --------------------------------------------------
$ws.Range("A3").Value = "Seeed XIAO (x3)"
$ws.Range("B3").Value = 25
$ws.Range("D3").Value = "Seeed XIAO (x3)"
$ws.Range("E3").Formula = "=B3"

# --- New component labels, entered in the order the author typed them
# (keeps the shared-string table ordering identical to the source edit) --
$ws.Range("A4").Value = "NRF24 SMD (x10)"
$ws.Range("D7").Value = "AMS1117"
$ws.Range("A7").Value = "AMS1117 (x10)"
$ws.Range("A8").Value = "Power Banks (x3)"
$ws.Range("A9").Value = "LEDs (x450)"
$ws.Range("A10").Value = "1k ohm (x100)"
$ws.Range("D4").Value = "NRF24 SMD (x3)"
$ws.Range("D8").Value = "Power Banks"
$ws.Range("D9").Value = "LEDs (x4)"
$ws.Range("D10").Value = "Resistors (x6)"

# --- Row 4: NRF24 SMD -----------------------------------------------------
$ws.Range("B4").Value = 10.88
$ws.Range("E4").Formula = "=(B4/10)*3"
# New thin left/right border on A4 (matches the new borderId in styles.xml)
$ws.Range("A4").Borders.Item(7).Weight = 2   # xlEdgeLeft, xlThin
$ws.Range("A4").Borders.Item(10).Weight = 2  # xlEdgeRight, xlThin

# --- Row 5: Buzzer --------------------------------------------------------
$ws.Range("A5").Value = "Buzzer (x10)"
$ws.Range("B5").Value = 7
$ws.Range("D5").Value = "Buzzer (x1)"
$ws.Range("E5").Formula = "=7/10"

# --- Row 6: Banana Jacks ---------------------------------------------------
$ws.Range("A6").Value = "Banana Jacks (x20)"
$ws.Range("B6").Value = 11
$ws.Range("D6").Value = "Banana Jacks (x6)"
$ws.Range("E6").Formula = "=(B6/20)*6"

# --- Row 7: AMS1117 ----------------------------------------------------
$ws.Range("B7").Value = 8
$ws.Range("E7").Formula = "=B7/10*3"

# --- Row 8: Power Banks -----------------------------------------------
$ws.Range("B8").Value = 26
$ws.Range("E8").Value = 26

# --- Row 9: LEDs (new row) ---------------------------------------------
$ws.Range("B9").Value = 13
$ws.Range("E9").Formula = "=B9/450*4"

# --- Row 10: 1k ohm resistors (new row) ---------------------------------
$ws.Range("B10").Value = 5
$ws.Range("E10").Formula = "=B10/100*6"

# --- Row 14: totals (moved down from row 10, now a gap of rows 12-13) ---
$ws.Range("A14").Value = "Total Order"
$ws.Range("B14").Formula = "=SUM(B3:B10)"
$ws.Range("D14").Value = "Per Unit"
$ws.Range("E14").Formula = "=SUM(E3:E8)"

# --- Selection / view state ------------------------------------------
$ws.Range("H26").Select()
